$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Login"
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Password"
$ws.Range("B2").Select()
